$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Sheet1 (quality_comparison) ---
# C1 -> top+bottom border (border 4)
$ws1.Range("C1").Style = "Normal"
$ws1.Range("C1").Borders.LineStyle = 1
$ws1.Range("C1").Borders.Weight = 2
$ws1.Range("C1").Borders.Item(7).LineStyle = -4142
$ws1.Range("C1").Borders.Item(10).LineStyle = -4142

# D1 -> top+bottom+right border (border 5)
$ws1.Range("D1").Style = "Normal"
$ws1.Range("D1").Borders.LineStyle = 1
$ws1.Range("D1").Borders.Weight = 2
$ws1.Range("D1").Borders.Item(7).LineStyle = -4142

# C2 text: fedcore -> approach
$ws1.Range("C2").Value = "approach"

# --- Sheet2 (computational_comparison) ---
# C1 -> top+bottom border (border 4)
$ws2.Range("C1").Style = "Normal"
$ws2.Range("C1").Borders.LineStyle = 1
$ws2.Range("C1").Borders.Weight = 2
$ws2.Range("C1").Borders.Item(7).LineStyle = -4142
$ws2.Range("C1").Borders.Item(10).LineStyle = -4142

# D1 -> top+bottom+right border (border 5)
$ws2.Range("D1").Style = "Normal"
$ws2.Range("D1").Borders.LineStyle = 1
$ws2.Range("D1").Borders.Weight = 2
$ws2.Range("D1").Borders.Item(7).LineStyle = -4142

# F1 -> top+bottom border (border 4)
$ws2.Range("F1").Style = "Normal"
$ws2.Range("F1").Borders.LineStyle = 1
$ws2.Range("F1").Borders.Weight = 2
$ws2.Range("F1").Borders.Item(7).LineStyle = -4142
$ws2.Range("F1").Borders.Item(10).LineStyle = -4142

# G1 -> top+bottom+right border (border 5)
$ws2.Range("G1").Style = "Normal"
$ws2.Range("G1").Borders.LineStyle = 1
$ws2.Range("G1").Borders.Weight = 2
$ws2.Range("G1").Borders.Item(7).LineStyle = -4142

# C2, F2 text: fedcore -> approach
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 cell removed entirely (was empty inlineStr)
$ws2.Range("G5").ClearContents()

Write-Host "All edits applied"
